# Setup LAb e soluzione
#
# 1) Remove the two empty paragraphs that precede the
#    "Request-Response Asynch One To One" paragraph.
# 2) Move the "_GoBack" bookmark from the end of the document
#    (after the "(selectors)" run) to the start of that same
#    "Request-Response Asynch One To One" paragraph (right after
#    its paragraph properties, before its first run). Because
#    bookmark names are unique, re-adding "_GoBack" at the new
#    location automatically relocates it away from its old spot.

$d = $word.ActiveDocument

# --- Step 1: delete the two blank paragraphs (originally paragraphs 4 & 5) ---
$blank1 = $d.Paragraphs.Item(4)
$d.Range($blank1.Range.Start, $blank1.Range.End).Delete()

$blank2 = $d.Paragraphs.Item(4)
$d.Range($blank2.Range.Start, $blank2.Range.End).Delete()

# --- Step 2: move the "_GoBack" bookmark to the target paragraph ---
$target = $d.Paragraphs.Item(4)
$start = $target.Range.Start
$bookmarkRange = $d.Range($start, $start)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
